$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog-Stories")

# Add the link reference in H6 first so it is appended to shared strings before the
# updated C6 text (matches shared string ordering produced by the original edit).
$ws.Range("H6").Value = "https://quangnguyennd.medium.com/git-flow-vs-github-flow-620c922b2cbd"

# Update the text in C6: drop "but simpler"
$ws.Range("C6").Value = "Based on github-flow"

# Update the selection on the active sheet to C7 (matches recorded selection change)
$ws.Activate()
$ws.Range("C7").Select()
